$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counter values on row 2 (MD state counters)
$ws.Range("B2").Value = 9483
$ws.Range("C2").Value = 2279
